$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 2376.2263
$ws.Cells.Item(17, 10).Value = 2409.173
$ws.Cells.Item(17, 12).Value = 7227.518999999999
$ws.Cells.Item(17, 14).Value = -7563.518999999999

# Row 129
$ws.Cells.Item(129, 8).Value = 839.2174
$ws.Cells.Item(129, 9).Value = 482.83334
$ws.Cells.Item(129, 11).Value = 1448.50002
$ws.Cells.Item(129, 13).Value = 3551.49998

# Row 135
$ws.Cells.Item(135, 8).Value = 31260032
$ws.Cells.Item(135, 9).Value = 1046.909
$ws.Cells.Item(135, 11).Value = 9422.181
$ws.Cells.Item(135, 13).Value = -6887.181

# Row 138
$ws.Cells.Item(138, 8).Value = 2589.6956
$ws.Cells.Item(138, 9).Value = 811.1539
$ws.Cells.Item(138, 10).Value = 3290.3333
$ws.Cells.Item(138, 11).Value = 2433.4617
$ws.Cells.Item(138, 12).Value = 9870.999899999999
$ws.Cells.Item(138, 13).Value = 2706.5383
$ws.Cells.Item(138, 14).Value = -20150.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Cells.Item(22, 8).Value = 1463
$ws.Cells.Item(22, 9).Value = 1463
$ws.Cells.Item(22, 11).Value = 1463
$ws.Cells.Item(22, 13).Value = -1164

# Row 32
$ws.Cells.Item(32, 8).Value = 18635.836
$ws.Cells.Item(32, 9).Value = 21499.354
$ws.Cells.Item(32, 11).Value = 21499.354
$ws.Cells.Item(32, 13).Value = -21212.354

# Row 45
$ws.Cells.Item(45, 8).Value = 2738.3235
$ws.Cells.Item(45, 9).Value = 1899.2142
$ws.Cells.Item(45, 10).Value = 3325.7
$ws.Cells.Item(45, 11).Value = 1899.2142
$ws.Cells.Item(45, 12).Value = 3325.7
$ws.Cells.Item(45, 13).Value = -1522.2142
$ws.Cells.Item(45, 14).Value = -4079.7

# Row 74
$ws.Cells.Item(74, 8).Value = 3938.7058
$ws.Cells.Item(74, 9).Value = 5144.8
$ws.Cells.Item(74, 11).Value = 5144.8
$ws.Cells.Item(74, 13).Value = -4270.8

# Row 77
$ws.Cells.Item(77, 8).Value = 3938.7058
$ws.Cells.Item(77, 9).Value = 5144.8
$ws.Cells.Item(77, 11).Value = 25724
$ws.Cells.Item(77, 13).Value = -21356

# Row 122
$ws.Cells.Item(122, 8).Value = 2336.9697
$ws.Cells.Item(122, 9).Value = 2220.28
$ws.Cells.Item(122, 10).Value = 2701.625
$ws.Cells.Item(122, 11).Value = 6660.84
$ws.Cells.Item(122, 12).Value = 8104.875
$ws.Cells.Item(122, 13).Value = -4210.84
$ws.Cells.Item(122, 14).Value = -13004.875

# Row 139
$ws.Cells.Item(139, 8).Value = 50624.645
$ws.Cells.Item(139, 10).Value = 50624.645
$ws.Cells.Item(139, 12).Value = 50624.645
$ws.Cells.Item(139, 14).Value = -60904.645

$ws = $wb.Worksheets.Item("BSM")
# Row 6
$ws.Cells.Item(6, 8).Value = 10443.667
$ws.Cells.Item(6, 10).Value = 10443.667
$ws.Cells.Item(6, 12).Value = 10443.667
$ws.Cells.Item(6, 14).Value = -10669.667

# Row 134
$ws.Cells.Item(134, 8).Value = 67176.56
$ws.Cells.Item(134, 9).Value = 67176.56
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 201529.68
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -198994.68
$ws.Cells.Item(134, 14).Value = $null

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 9847.137000000001
$ws.Cells.Item(31, 9).Value = 20105.611
$ws.Cells.Item(31, 10).Value = 2745.1155
$ws.Cells.Item(31, 11).Value = 20105.611
$ws.Cells.Item(31, 12).Value = 2745.1155
$ws.Cells.Item(31, 13).Value = -19810.611
$ws.Cells.Item(31, 14).Value = -3335.1155

# Row 34
$ws.Cells.Item(34, 8).Value = 9847.137000000001
$ws.Cells.Item(34, 9).Value = 20105.611
$ws.Cells.Item(34, 10).Value = 2745.1155
$ws.Cells.Item(34, 11).Value = 20105.611
$ws.Cells.Item(34, 12).Value = 2745.1155
$ws.Cells.Item(34, 13).Value = -19903.611
$ws.Cells.Item(34, 14).Value = -3149.1155

# Row 110
$ws.Cells.Item(110, 8).Value = 50000
$ws.Cells.Item(110, 10).Value = 50000
$ws.Cells.Item(110, 12).Value = 50000
$ws.Cells.Item(110, 14).Value = -58180

# Row 119
$ws.Cells.Item(119, 8).Value = 50000
$ws.Cells.Item(119, 10).Value = 50000
$ws.Cells.Item(119, 12).Value = 50000
$ws.Cells.Item(119, 14).Value = -59676

# Row 120
$ws.Cells.Item(120, 8).Value = 24959.7
$ws.Cells.Item(120, 9).Value = 11715
$ws.Cells.Item(120, 10).Value = 30636
$ws.Cells.Item(120, 11).Value = 11715
$ws.Cells.Item(120, 12).Value = 30636
$ws.Cells.Item(120, 13).Value = -8086
$ws.Cells.Item(120, 14).Value = -37894

# Row 121
$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 9).Value = 0
$ws.Cells.Item(121, 11).Value = 0
$ws.Cells.Item(121, 13).Value = $null

$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Cells.Item(6, 8).Value = 80.57143000000001
$ws.Cells.Item(6, 9).Value = 63.53846
$ws.Cells.Item(6, 11).Value = 190.61538
$ws.Cells.Item(6, 13).Value = -77.61538000000002

# Row 16
$ws.Cells.Item(16, 8).Value = 500
$ws.Cells.Item(16, 10).Value = 500
$ws.Cells.Item(16, 12).Value = 1500
$ws.Cells.Item(16, 14).Value = -1846

# Row 68
$ws.Cells.Item(68, 8).Value = 3504.634
$ws.Cells.Item(68, 9).Value = 966.3333
$ws.Cells.Item(68, 10).Value = 3705.0264
$ws.Cells.Item(68, 11).Value = 2898.9999
$ws.Cells.Item(68, 12).Value = 11115.0792
$ws.Cells.Item(68, 13).Value = -2087.9999
$ws.Cells.Item(68, 14).Value = -12737.0792

# Row 71
$ws.Cells.Item(71, 8).Value = 3504.634
$ws.Cells.Item(71, 9).Value = 966.3333
$ws.Cells.Item(71, 10).Value = 3705.0264
$ws.Cells.Item(71, 11).Value = 8696.9997
$ws.Cells.Item(71, 12).Value = 33345.2376
$ws.Cells.Item(71, 13).Value = -4640.9997
$ws.Cells.Item(71, 14).Value = -41457.2376

# Row 107
$ws.Cells.Item(107, 8).Value = 4273.3
$ws.Cells.Item(107, 9).Value = 33633.332
$ws.Cells.Item(107, 10).Value = 1011.0741
$ws.Cells.Item(107, 11).Value = 100899.996
$ws.Cells.Item(107, 12).Value = 3033.2223
$ws.Cells.Item(107, 13).Value = -98979.99600000001
$ws.Cells.Item(107, 14).Value = -6873.2223

# Row 113
$ws.Cells.Item(113, 8).Value = 13030
$ws.Cells.Item(113, 9).Value = 25387.5
$ws.Cells.Item(113, 10).Value = 672.5
$ws.Cells.Item(113, 11).Value = 76162.5
$ws.Cells.Item(113, 12).Value = 2017.5
$ws.Cells.Item(113, 13).Value = -73992.5
$ws.Cells.Item(113, 14).Value = -6357.5

# Row 131
$ws.Cells.Item(131, 8).Value = 842.13
$ws.Cells.Item(131, 9).Value = 625
$ws.Cells.Item(131, 10).Value = 851.17706
$ws.Cells.Item(131, 11).Value = 1875
$ws.Cells.Item(131, 12).Value = 2553.53118
$ws.Cells.Item(131, 13).Value = 3165
$ws.Cells.Item(131, 14).Value = -12633.53118

# Row 132
$ws.Cells.Item(132, 8).Value = 750.2105
$ws.Cells.Item(132, 10).Value = 812.5
$ws.Cells.Item(132, 12).Value = 7312.5
$ws.Cells.Item(132, 14).Value = -12372.5

# Row 140
$ws.Cells.Item(140, 8).Value = 1908.3125
$ws.Cells.Item(140, 9).Value = 1508.0769
$ws.Cells.Item(140, 11).Value = 4524.2307
$ws.Cells.Item(140, 13).Value = 655.7692999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 14).Value = $null

# Row 24
$ws.Cells.Item(24, 8).Value = 37333.332
$ws.Cells.Item(24, 10).Value = 6000
$ws.Cells.Item(24, 12).Value = 6000
$ws.Cells.Item(24, 14).Value = -6346

# Row 32
$ws.Cells.Item(32, 8).Value = 23500
$ws.Cells.Item(32, 10).Value = 23500
$ws.Cells.Item(32, 12).Value = 23500
$ws.Cells.Item(32, 14).Value = -24092

# Row 113
$ws.Cells.Item(113, 8).Value = 3800.9
$ws.Cells.Item(113, 9).Value = 2168.3333
$ws.Cells.Item(113, 11).Value = 2168.3333
$ws.Cells.Item(113, 13).Value = 1.666700000000219

# Row 122
$ws.Cells.Item(122, 8).Value = 3458.8667
$ws.Cells.Item(122, 9).Value = 2657.75
$ws.Cells.Item(122, 11).Value = 7973.25
$ws.Cells.Item(122, 13).Value = -5523.25

# Row 132
$ws.Cells.Item(132, 8).Value = 97696.625
$ws.Cells.Item(132, 9).Value = 80934.30499999999
$ws.Cells.Item(132, 10).Value = 170333.33
$ws.Cells.Item(132, 11).Value = 242802.915
$ws.Cells.Item(132, 12).Value = 510999.99
$ws.Cells.Item(132, 13).Value = -240272.915
$ws.Cells.Item(132, 14).Value = -516059.99

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 2080.3333
$ws.Cells.Item(22, 9).Value = 2320
$ws.Cells.Item(22, 11).Value = 2320
$ws.Cells.Item(22, 13).Value = -2025

# Row 27
$ws.Cells.Item(27, 8).Value = 2080.3333
$ws.Cells.Item(27, 9).Value = 2320
$ws.Cells.Item(27, 11).Value = 2320
$ws.Cells.Item(27, 13).Value = -2213

# Row 40
$ws.Cells.Item(40, 8).Value = 73185.06
$ws.Cells.Item(40, 9).Value = 88654.69500000001
$ws.Cells.Item(40, 10).Value = 6150
$ws.Cells.Item(40, 11).Value = 88654.69500000001
$ws.Cells.Item(40, 12).Value = 6150
$ws.Cells.Item(40, 13).Value = -88518.69500000001
$ws.Cells.Item(40, 14).Value = -6422

# Row 93
$ws.Cells.Item(93, 8).Value = 1918.8948
$ws.Cells.Item(93, 9).Value = 2051.8125
$ws.Cells.Item(93, 10).Value = 1210
$ws.Cells.Item(93, 11).Value = 2051.8125
$ws.Cells.Item(93, 12).Value = 1210
$ws.Cells.Item(93, 13).Value = -803.8125
$ws.Cells.Item(93, 14).Value = -3706

# Row 122
$ws.Cells.Item(122, 8).Value = 3283.875
$ws.Cells.Item(122, 9).Value = 2581.889
$ws.Cells.Item(122, 11).Value = 7745.667
$ws.Cells.Item(122, 13).Value = -5295.667

# Row 136
$ws.Cells.Item(136, 8).Value = 2099.476
$ws.Cells.Item(136, 9).Value = 1592.3334
$ws.Cells.Item(136, 10).Value = 3367.3333
$ws.Cells.Item(136, 11).Value = 4777.0002
$ws.Cells.Item(136, 12).Value = 10101.9999
$ws.Cells.Item(136, 13).Value = -2227.0002
$ws.Cells.Item(136, 14).Value = -15201.9999

# Row 141
$ws.Cells.Item(141, 8).Value = 58357.5
$ws.Cells.Item(141, 10).Value = 58357.5
$ws.Cells.Item(141, 12).Value = 58357.5
$ws.Cells.Item(141, 14).Value = -68717.5

$ws = $wb.Worksheets.Item("WVR")
# Row 22
$ws.Cells.Item(22, 8).Value = 2000
$ws.Cells.Item(22, 9).Value = 2000
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 2000
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -1707
$ws.Cells.Item(22, 14).Value = $null

# Row 31
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 13).Value = $null

# Row 55
$ws.Cells.Item(55, 8).Value = 10000
$ws.Cells.Item(55, 9).Value = 1000
$ws.Cells.Item(55, 10).Value = 11800
$ws.Cells.Item(55, 11).Value = 1000
$ws.Cells.Item(55, 12).Value = 11800
$ws.Cells.Item(55, 13).Value = -723
$ws.Cells.Item(55, 14).Value = -12354

# Row 86
$ws.Cells.Item(86, 8).Value = 9770
$ws.Cells.Item(86, 10).Value = 9770
$ws.Cells.Item(86, 12).Value = 9770
$ws.Cells.Item(86, 14).Value = -12016

# Row 89
$ws.Cells.Item(89, 8).Value = 9770
$ws.Cells.Item(89, 10).Value = 9770
$ws.Cells.Item(89, 12).Value = 48850
$ws.Cells.Item(89, 14).Value = -60082

# Row 109
$ws.Cells.Item(109, 8).Value = 27330
$ws.Cells.Item(109, 10).Value = 27330
$ws.Cells.Item(109, 12).Value = 27330
$ws.Cells.Item(109, 14).Value = -30104

# Row 136
$ws.Cells.Item(136, 8).Value = 1231.32
$ws.Cells.Item(136, 9).Value = 892.7273
$ws.Cells.Item(136, 10).Value = 1497.3572
$ws.Cells.Item(136, 11).Value = 2678.1819
$ws.Cells.Item(136, 12).Value = 4492.071599999999
$ws.Cells.Item(136, 13).Value = -128.1819
$ws.Cells.Item(136, 14).Value = -9592.071599999999
